$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(7).Delete()
$ws.Columns.Item(11).Delete()
Write-Host ("C1 before insert: " + $ws.Range("C1").Value2)
$ws.Columns.Item(3).Insert()
Write-Host ("C1 after insert: [" + $ws.Range("C1").Value2 + "]")
Write-Host ("D1 after insert: " + $ws.Range("D1").Value2)
